$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for team record: Wins, Losses, Ties
# Copy formatting (bold header style) from the existing last header cell (AC1)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record values for every player row (2-52)
$ws.Range("AD2:AD52").Value = 101
$ws.Range("AE2:AE52").Value = 61
$ws.Range("AF2:AF52").Value = 0
